$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the hyperlinks that used to sit on B2 and C2 (only B5's mailto
#    hyperlink should remain afterwards).
# ---------------------------------------------------------------------------
$found = $true
while ($found) {
    $found = $false
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($addr -eq '$B$2' -or $addr -eq '$C$2') {
            $hl.Delete()
            $found = $true
            break
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Update row 2 - jagdish.d / User@12345 / Baramati login / BMC
# ---------------------------------------------------------------------------
$ws.Range("B2").Value2 = "User@12345"
$ws.Range("C2").Value2 = "http://testbaramatimc.ptaxcollection.com:8080/Pages/Login.aspx"
$ws.Range("D2").Value2 = "BMC"

# ---------------------------------------------------------------------------
# 3. Rows 10-13: replace the "KM" property block with the new "BMC" data.
#    (values only change here - content was previously Cash/Cheque-clear/
#    Card/advance and becomes Cash+advance/Cheque-bounce/Cheque-Clear/Card)
# ---------------------------------------------------------------------------
$ws.Range("A10").Value2 = "BMC"
$ws.Range("B10").Value2 = 4
$ws.Range("C10").Value2 = "18"
$ws.Range("D10").Value2 = "Cash+advance"
$ws.Range("E10").Value2 = "node1"

$ws.Range("A11").Value2 = "BMC"
$ws.Range("B11").Value2 = 4
$ws.Range("C11").Value2 = "17"
$ws.Range("D11").Value2 = "Cheque-bounce"
$ws.Range("E11").Value2 = "node2"

$ws.Range("A12").Value2 = "BMC"
$ws.Range("B12").Value2 = 4
$ws.Range("C12").Value2 = "20"
$ws.Range("D12").Value2 = "Cheque-Clear"
$ws.Range("E12").Value2 = "node3"

# Row 13 used to be a normal (non-bottom-border) row; it now becomes the
# last row of the shrunk table, matching the borders/height used by rows
# 11-12 (thick bottom border).
$ws.Range("A13").Value2 = "BMC"
$ws.Range("B13").Value2 = 4
$ws.Range("C13").Value2 = "99"
$ws.Range("D13").Value2 = "Card"
$ws.Range("E13").Value2 = "node4"

$ws.Range("A12:E12").Copy() | Out-Null
$ws.Range("A13:E13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# B11 and B12 previously carried a bottom/medium border left over from the
# old layout (rows 11/12 used to be mid-table); in the new 4-row table only
# the whole-row thick-bottom border (set above) should remain, so clear the
# stray per-cell borders that used to be there by re-applying the plain
# (border-less) look used by the A column.
$ws.Range("A10").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null
$ws.Range("B12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# C11/C12 already carry the right+bottom medium border combination that C13
# now needs too (copied above); re-apply it to C13 alone so the value cell
# uses exactly the same border pattern as its neighbours.
$ws.Range("C12").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Delete rows 14-24: the second "BMC"/duplicate table plus the trailing
#    blank rows are gone entirely - the sheet now ends at row 13.
# ---------------------------------------------------------------------------
$ws.Rows("14:24").Delete() | Out-Null

# ---------------------------------------------------------------------------
# 5. Sheet view bookkeeping - zoom + active selection moved to reflect the
#    smaller sheet.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 145
$ws.Range("C14").Select() | Out-Null

$wb.Save()
